# Apply corrections to "Dados das NFs" sheet:
#  1. Remove the two duplicate/erroneous rows (old rows 27 and 28), shifting
#     the remaining rows up.
#  2. Reformat every 14-digit CNPJ value in columns C and D (rows 2-33)
#     from plain digits to the standard NN.NNN.NNN/NNNN-NN mask.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two bad rows (329557 / 336834) -------------------------
$ws.Rows("27:28").Delete()

# --- 2. Reformat CNPJ values in columns C and D (rows 2-33) ---------------
function Format-Cnpj([string]$digits) {
    $p1 = $digits.Substring(0,2)
    $p2 = $digits.Substring(2,3)
    $p3 = $digits.Substring(5,3)
    $p4 = $digits.Substring(8,4)
    $p5 = $digits.Substring(12,2)
    return "{0}.{1}.{2}/{3}-{4}" -f $p1,$p2,$p3,$p4,$p5
}

for ($row = 2; $row -le 33; $row++) {
    foreach ($col in "C", "D") {
        $cell = $ws.Range($col + $row)
        $raw = [string]$cell.Value()
        if ($raw -ne $null -and $raw.Length -eq 14) {
            $cell.Value = Format-Cnpj $raw
        }
    }
}
